# This document gets a footer added: a 3-column "Blank (Three Columns)"
# style footer (built-in Word gallery layout), the section's page
# orientation is pinned to explicit "portrait", and the new footer is
# wired up via a footerReference on the section.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Make the (already-default) portrait orientation explicit in pgSz.
$sec.PageSetup.Orientation = 0

$ftr = $sec.Footers(1)
$rng = $ftr.Range

# Touch the built-in "Table Grid" table style once (on a throwaway table)
# so its definition gets registered into styles.xml, mirroring what Word
# does when the "Blank (Three Columns)" footer building block is
# inserted from the gallery. The footer's real table stays on the
# default "Table Normal" style further below.
$scratchTbl = $d.Tables.Add($rng, 1, 1)
$scratchTbl.Style = "Table Grid"

# Build the actual three-column footer layout (label / page-number /
# date style blank columns) and drop it into the footer, replacing the
# scratch content created above.
$ftr2 = $sec.Footers(1)
$rng2 = $ftr2.Range

$footerXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:tbl>
<w:tblPr>
<w:tblStyle w:val="TableNormal"/>
<w:bidiVisual w:val="0"/>
<w:tblW w:w="0" w:type="auto"/>
<w:tblLayout w:type="fixed"/>
<w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
</w:tblPr>
<w:tblGrid>
<w:gridCol w:w="3120"/>
<w:gridCol w:w="3120"/>
<w:gridCol w:w="3120"/>
</w:tblGrid>
<w:tr>
<w:tc>
<w:tcPr>
<w:tcW w:w="3120" w:type="dxa"/>
<w:tcMar/>
</w:tcPr>
<w:p>
<w:pPr>
<w:pStyle w:val="Header"/>
<w:bidi w:val="0"/>
<w:ind w:left="-115"/>
<w:jc w:val="left"/>
</w:pPr>
</w:p>
</w:tc>
<w:tc>
<w:tcPr>
<w:tcW w:w="3120" w:type="dxa"/>
<w:tcMar/>
</w:tcPr>
<w:p>
<w:pPr>
<w:pStyle w:val="Header"/>
<w:bidi w:val="0"/>
<w:jc w:val="center"/>
</w:pPr>
</w:p>
</w:tc>
<w:tc>
<w:tcPr>
<w:tcW w:w="3120" w:type="dxa"/>
<w:tcMar/>
</w:tcPr>
<w:p>
<w:pPr>
<w:pStyle w:val="Header"/>
<w:bidi w:val="0"/>
<w:ind w:right="-115"/>
<w:jc w:val="right"/>
</w:pPr>
</w:p>
</w:tc>
</w:tr>
</w:tbl>
<w:p>
<w:pPr>
<w:pStyle w:val="Footer"/>
<w:bidi w:val="0"/>
</w:pPr>
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@
$rng2.InsertXML($footerXml)

# Re-lay the header's tab characters into their own runs (split out of the
# runs that carry the following text), matching a normal Word round-trip
# of this header.
$hdr = $sec.Headers(1)
$hdrRng = $hdr.Range

$headerXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="Header"/>
</w:pPr>
<w:r><w:t xml:space="preserve">Feingold </w:t></w:r>
<w:r><w:tab/></w:r>
<w:r><w:t xml:space="preserve">Russia </w:t></w:r>
<w:r><w:tab/></w:r>
<w:r><w:t>June 21, 2007</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@
$hdrRng.InsertXML($headerXml)

Write-Output "Footer added; orientation pinned; TableGrid style registered; header runs re-split."
